$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new member row (row 12) for Scott Anderson, riding of
# Vernon-Lake Country-Monashee. He has no phone number on file, so the
# Phone column is left blank (workaround for people with no phone or
# email) rather than omitted.
$ws.Range("A12").Value = "Anderson, Scott"
$ws.Range("B12").Value = "https://www.ourcommons.ca/Members/en/scott-anderson(89259)"
$ws.Range("C12").Value = "Conservative"
$ws.Range("D12").Value = "Vernon—Lake Country—Monashee"
$ws.Range("E12").Value = "https://www.ourcommons.ca/Members/en/constituencies/vernon-lake-country-monashee(1280)"
$ws.Range("F12").Value = "British Columbia"

# Phone column: no phone on file. Assigning a plain "" removes the cell
# entirely instead of leaving a real (blank) text cell in place, so force
# a quote-prefixed empty string first (guarantees an actual empty-text
# cell gets written) and then clear the quote-prefix formatting it
# introduces so the cell is left with the default style.
$ws.Range("G12").Value = "'"
$ws.Range("G12").ClearFormats()

$ws.Range("H12").Value = "scott.anderson@parl.gc.ca"
